$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.105.01"
$ws.Range("E2").Value = "  +3.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.723.84"
$ws.Range("E3").Value = "  +2.60%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.93"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.33"
$ws.Range("E8").Value = "  +13.73%  "

$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0633"
$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.967.90"
$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.726.09"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("E14").Value = "  +2.94%  "

$ws.Range("E15").Value = "  +4.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.55"
$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.069.52"
$ws.Range("E17").Value = "  +3.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.87"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0754"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.66"
$ws.Range("E23").Value = "  +1.93%  "

$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.10"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("E26").Value = "  +3.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.69"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  +2.02%  "

$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("E32").Value = "  +1.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.499.18"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.27"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  +2.19%  "

$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.79"
$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.81"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.872.47"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.803"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.76"
$ws.Range("E47").Value = "  +10.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "90.81"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0113"
$ws.Range("E49").Value = "  +4.84%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.27"
$ws.Range("E50").Value = "  +1.82%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  +0.34%  "
